$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# NOTE: all insertions below are performed from the bottom of the document
# upward so that paragraph indices captured for the earlier (higher up)
# edits remain valid while we work.
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# Step D: after the last numbered item ("There would need to be enough
# resources ...", paragraph 27) and before the trailing tab-only paragraph,
# insert a new numbered list paragraph (numId 5, inherited) with the new
# answer text, a bookmark ("_GoBack") right after its lead-in phrase, and a
# trailing indented (but otherwise plain) blank paragraph.
# ---------------------------------------------------------------------------
$pLastItem = $d.Paragraphs(27)
$rLastItem = $pLastItem.Range
$rLastItem.Collapse(0)
$rLastItem.InsertParagraphAfter()

$pNewList = $d.Paragraphs(28)
$pNewList.Range.InsertAfter("In this case, many additional variables were needed. First of all, this is a bit field that determines the type of operation. Then each specific type of operation requires its own set of flags. For mathematic operation they are type of operation source / destination address, operation width (8/16 bits). For memory operation - direction (store/load), register and memory adressing mode. And for branch operation -type of branching.")

$bmParaStart = $d.Paragraphs(28).Range.Start
$bmPoint = $bmParaStart + 14
$bmRange = $d.Range($bmPoint, $bmPoint)
$d.Bookmarks.Add("_GoBack", $bmRange)

$pNewList2 = $d.Paragraphs(28)
$rNewList2 = $pNewList2.Range
$rNewList2.Collapse(0)
$rNewList2.InsertParagraphAfter()
$pBlankD = $d.Paragraphs(29)
$pBlankD.Range.ParagraphFormat.Style = "Normal"
$pBlankD.Range.ParagraphFormat.LeftIndent = 18

# ---------------------------------------------------------------------------
# Step C: before "4. If executeInstruction() ..." (paragraph 25), insert a
# new blank paragraph indented the same way as the others.
# ---------------------------------------------------------------------------
$p24 = $d.Paragraphs(24)
$r24 = $p24.Range
$r24.Collapse(0)
$r24.InsertParagraphAfter()
$pBlankC = $d.Paragraphs(25)
$pBlankC.Range.ParagraphFormat.Style = "Normal"
$pBlankC.Range.ParagraphFormat.LeftIndent = 18

# ---------------------------------------------------------------------------
# Step B: after the "121 unused opcodes ..." answer (paragraph 17) and
# before "3. What would we need ..." insert a new blank indented paragraph.
# ---------------------------------------------------------------------------
$p17 = $d.Paragraphs(17)
$r17 = $p17.Range
$r17.Collapse(0)
$r17.InsertParagraphAfter()
$pBlankB = $d.Paragraphs(18)
$pBlankB.Range.ParagraphFormat.Style = "Normal"
$pBlankB.Range.ParagraphFormat.LeftIndent = 18

# ---------------------------------------------------------------------------
# Step A: after "STOR ACC, [0x0000]" (paragraph 13) and before "2. Of the
# 256 ..." insert a plain blank paragraph.
# ---------------------------------------------------------------------------
$p13 = $d.Paragraphs(13)
$r13 = $p13.Range
$r13.Collapse(0)
$r13.InsertParagraphAfter()

Write-Output "edits applied"
